# Auto-generated Excel COM-interop script
# Applies updated market price / profit values to the Leve profit sheets
# (values refreshed by the scheduled Sheets runner).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 330.4
$ws.Range("I31").Value = 330.4
$ws.Range("K31").Value = 991.1999999999999
$ws.Range("M31").Value = -761.1999999999999
$ws.Range("H113").Value = 3292
$ws.Range("I113").Value = 3022.6667
$ws.Range("K113").Value = 3022.6667
$ws.Range("M113").Value = 231.3332999999998
$ws.Range("H121").Value = 1281.7646
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1281.7646
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 3845.2938
$ws.Range("M121").ClearContents() | Out-Null
$ws.Range("N121").Value = -7339.293799999999
$ws.Range("H132").Value = 1326109.5
$ws.Range("I132").Value = 1685.0667
$ws.Range("J132").Value = 7002214
$ws.Range("K132").Value = 5055.2001
$ws.Range("L132").Value = 21006642
$ws.Range("M132").Value = -2525.2001
$ws.Range("N132").Value = -21011702
$ws.Range("H137").Value = 4351784.5
$ws.Range("I137").Value = 6669889.5
$ws.Range("J137").Value = 5337.25
$ws.Range("K137").Value = 20009668.5
$ws.Range("L137").Value = 16011.75
$ws.Range("M137").Value = -20007118.5
$ws.Range("N137").Value = -21111.75
$ws.Range("H138").Value = 2851527.8
$ws.Range("I138").Value = 239748.69
$ws.Range("J138").Value = 8336263.5
$ws.Range("K138").Value = 719246.0700000001
$ws.Range("L138").Value = 25008790.5
$ws.Range("M138").Value = -714106.0700000001
$ws.Range("N138").Value = -25019070.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 4704.2856
$ws.Range("I21").Value = 3821.6667
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 3821.6667
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -3447.6667
$ws.Range("N21").Value = -10748
$ws.Range("H32").Value = 839.96
$ws.Range("I32").Value = 770.4888999999999
$ws.Range("J32").Value = 1465.2
$ws.Range("K32").Value = 770.4888999999999
$ws.Range("L32").Value = 1465.2
$ws.Range("M32").Value = -483.4888999999999
$ws.Range("N32").Value = -2039.2
$ws.Range("H74").Value = 5480294.5
$ws.Range("I74").Value = 6606679.5
$ws.Range("K74").Value = 6606679.5
$ws.Range("M74").Value = -6605805.5
$ws.Range("H77").Value = 5480294.5
$ws.Range("I77").Value = 6606679.5
$ws.Range("K77").Value = 33033397.5
$ws.Range("M77").Value = -33029029.5
$ws.Range("H102").Value = 10205966
$ws.Range("I102").Value = 12988593
$ws.Range("K102").Value = 12988593
$ws.Range("M102").Value = -12986971
$ws.Range("H122").Value = 2925666
$ws.Range("I122").Value = 1761.5625
$ws.Range("J122").Value = 18519824
$ws.Range("K122").Value = 5284.6875
$ws.Range("L122").Value = 55559472
$ws.Range("M122").Value = -2834.6875
$ws.Range("N122").Value = -55564372

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4254.963
$ws.Range("I31").Value = 2448.6667
$ws.Range("J31").Value = 5700
$ws.Range("K31").Value = 2448.6667
$ws.Range("L31").Value = 5700
$ws.Range("M31").Value = -2153.6667
$ws.Range("N31").Value = -6290
$ws.Range("H33").Value = 6000
$ws.Range("I33").Value = 6000
$ws.Range("K33").Value = 6000
$ws.Range("M33").Value = -5621
$ws.Range("H34").Value = 4254.963
$ws.Range("I34").Value = 2448.6667
$ws.Range("J34").Value = 5700
$ws.Range("K34").Value = 2448.6667
$ws.Range("L34").Value = 5700
$ws.Range("M34").Value = -2246.6667
$ws.Range("N34").Value = -6104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15152461
$ws.Range("I131").Value = 90909480
$ws.Range("K131").Value = 272728440
$ws.Range("M131").Value = -272723400

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3532.7585
$ws.Range("I80").Value = 2890
$ws.Range("J80").Value = 3580.3704
$ws.Range("K80").Value = 2890
$ws.Range("L80").Value = 3580.3704
$ws.Range("M80").Value = -1892
$ws.Range("N80").Value = -5576.3704
$ws.Range("H83").Value = 3532.7585
$ws.Range("I83").Value = 2890
$ws.Range("J83").Value = 3580.3704
$ws.Range("K83").Value = 14450
$ws.Range("L83").Value = 17901.852
$ws.Range("M83").Value = -9458
$ws.Range("N83").Value = -27885.852
$ws.Range("H102").Value = 1424.5714
$ws.Range("I102").Value = 976.6667
$ws.Range("J102").Value = 1760.5
$ws.Range("K102").Value = 976.6667
$ws.Range("L102").Value = 1760.5
$ws.Range("M102").Value = 645.3333
$ws.Range("N102").Value = -5004.5
$ws.Range("H136").Value = 6632
$ws.Range("J136").Value = 6632
$ws.Range("L136").Value = 19896
$ws.Range("N136").Value = -24996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3550.7144
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 5618.3335
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 5618.3335
$ws.Range("M16").Value = -1830
$ws.Range("N16").Value = -5958.3335
$ws.Range("H32").Value = 1895
$ws.Range("I32").Value = 1895
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1895
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1578
$ws.Range("N32").ClearContents() | Out-Null
$ws.Range("H68").Value = 1668.2174
$ws.Range("I68").Value = 1636.619
$ws.Range("K68").Value = 1636.619
$ws.Range("M68").Value = -887.6189999999999
$ws.Range("H71").Value = 1668.2174
$ws.Range("I71").Value = 1636.619
$ws.Range("K71").Value = 8183.094999999999
$ws.Range("M71").Value = -4439.094999999999
$ws.Range("H82").Value = 6600
$ws.Range("I82").Value = 1350
$ws.Range("J82").Value = 8350
$ws.Range("K82").Value = 1350
$ws.Range("L82").Value = 8350
$ws.Range("M82").Value = -989
$ws.Range("N82").Value = -9072
$ws.Range("H85").Value = 6600
$ws.Range("I85").Value = 1350
$ws.Range("J85").Value = 8350
$ws.Range("K85").Value = 1350
$ws.Range("L85").Value = 8350
$ws.Range("M85").Value = -102
$ws.Range("N85").Value = -10846
$ws.Range("H136").Value = 84056.92
$ws.Range("I136").Value = 49281.74
$ws.Range("K136").Value = 147845.22
$ws.Range("M136").Value = -145295.22

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2086.6458
$ws.Range("I122").Value = 1815.963
$ws.Range("J122").Value = 2434.6667
$ws.Range("K122").Value = 5447.889
$ws.Range("L122").Value = 7304.000100000001
$ws.Range("M122").Value = -2997.889
$ws.Range("N122").Value = -12204.0001
$ws.Range("H136").Value = 39117.605
$ws.Range("I136").Value = 31118.273
$ws.Range("J136").Value = 52316.5
$ws.Range("K136").Value = 93354.819
$ws.Range("L136").Value = 156949.5
$ws.Range("M136").Value = -90804.819
$ws.Range("N136").Value = -162049.5

